$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A ("Urutan")
$ws.Columns("A").Delete()

# Insert a new column after "Nama" (now column C), before "Jenis Kelamin" (now column D)
$ws.Columns("D").Insert()

# Set header & style for the newly inserted column D (bold, like the other header cells)
$ws.Range("D1").Value = "Tanggal Masuk Admedika"
$ws.Range("D1").Font.Bold = $true

# Re-select column A like in the target file
$ws.Range("A1:A1048576").Select()
